$d = $word.ActiveDocument

# 1) Objectives paragraph (English, italic) - split "...degree. " / "2 - Integration..."
$find1 = "1 - Consolidation and application of knowledge acquired in each of the specific areas of the Chemical Engineering degree. 2 - Integration of knowledge of Chemical Engineering"
$repl1 = "1 - Consolidation and application of knowledge acquired in each of the specific areas of the Chemical Engineering degree. ^l2 - Integration of knowledge of Chemical Engineering"
$d.Content.Find.Execute($find1, $false, $false, $false, $false, $false, $true, 1, $false, $repl1, 2)

# 2) Programa paragraph (Portuguese) - split into four segments at the numbered items
$find2 = "1 - Diagramas para estudos de processos químicos: diagramas de bloco; Fluxogramas de processo (PFD); Fluxogramas de instrumentação e tubulação (P&ID).2  Estrutura e síntese de processos químicos industriais: Hierarquia no planejamento de processos; Etapa 1- Descontínuo ou contínuo; Etapa 2 - Estrutura de entrada/saída de processo; Etapa 3- Estrutura de reciclo; 3  Análise de desempenho de processos químicos: Modelo de entrada e saída; Ferramentas para a avaliação de processos.4  Estudo de planta química industrial."
$repl2 = "1 - Diagramas para estudos de processos químicos: diagramas de bloco; Fluxogramas de processo (PFD); Fluxogramas de instrumentação e tubulação (P&ID).^l2  Estrutura e síntese de processos químicos industriais: Hierarquia no planejamento de processos; Etapa 1- Descontínuo ou contínuo; Etapa 2 - Estrutura de entrada/saída de processo; Etapa 3- Estrutura de reciclo; ^l3  Análise de desempenho de processos químicos: Modelo de entrada e saída; Ferramentas para a avaliação de processos.^l4  Estudo de planta química industrial."
$d.Content.Find.Execute($find2, $false, $false, $false, $false, $false, $true, 1, $false, $repl2, 2)

# 3) Programa paragraph (English, italic) - split into four segments at the numbered items
$find3 = "1 - Diagrams for Understanding Chemical Processes: Block Flow Diagrams; Process Flow Diagram (PFD); Piping and Instrumentation Diagram (P&ID).2 - The Structure and Synthesis of Process Flow Diagrams:  Hierarchy of Process Design; Step 1 - Batch versus Continuous Process; Step 2 - The Input/Output Structure of the Process; Step 3 - The Recycle Structure of the Process3 - Analysis of process performance: Process Input/Output Models; Tools for evaluating process performance.4 - Industrial chemical plant study."
$repl3 = "1 - Diagrams for Understanding Chemical Processes: Block Flow Diagrams; Process Flow Diagram (PFD); Piping and Instrumentation Diagram (P&ID).^l2 - The Structure and Synthesis of Process Flow Diagrams:  Hierarchy of Process Design; Step 1 - Batch versus Continuous Process; Step 2 - The Input/Output Structure of the Process; Step 3 - The Recycle Structure of the Process^l3 - Analysis of process performance: Process Input/Output Models; Tools for evaluating process performance.^l4 - Industrial chemical plant study."
$d.Content.Find.Execute($find3, $false, $false, $false, $false, $false, $true, 1, $false, $repl3, 2)
